$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "43.064.58"
$ws.Cells.Item(2, 5).Value = "  +1.39%  "
$ws.Cells.Item(3, 4).Value = "2.355.28"
$ws.Cells.Item(3, 5).Value = "  +5.91%  "
$ws.Cells.Item(4, 5).Value = "  -0.31%  "
$ws.Cells.Item(5, 4).Value = "'306.84"
$ws.Cells.Item(5, 5).Value = "  +2.73%  "
$ws.Cells.Item(6, 4).Value = "'107.06"
$ws.Cells.Item(6, 5).Value = "  -3.41%  "
$ws.Cells.Item(7, 5).Value = "  +1.10%  "
$ws.Cells.Item(8, 5).Value = "  -0.13%  "
$ws.Cells.Item(9, 4).Value = "'0.636"
$ws.Cells.Item(9, 5).Value = "  +4.27%  "
$ws.Cells.Item(10, 4).Value = "'42.66"
$ws.Cells.Item(10, 5).Value = "  -5.33%  "
$ws.Cells.Item(11, 4).Value = "'0.0936"
$ws.Cells.Item(11, 5).Value = "  +1.15%  "
$ws.Cells.Item(12, 4).Value = "'8.94"
$ws.Cells.Item(12, 5).Value = "  +1.13%  "
$ws.Cells.Item(13, 4).Value = "'1.05"
$ws.Cells.Item(13, 5).Value = "  +8.82%  "
$ws.Cells.Item(14, 5).Value = "  +1.14%  "
$ws.Cells.Item(15, 4).Value = "'16.44"
$ws.Cells.Item(15, 5).Value = "  +8.71%  "
$ws.Cells.Item(16, 4).Value = "2.714.58"
$ws.Cells.Item(16, 5).Value = "  +6.17%  "
$ws.Cells.Item(17, 4).Value = "2.441.21"
$ws.Cells.Item(17, 5).Value = "  +9.30%  "
$ws.Cells.Item(18, 4).Value = "43.071.34"
$ws.Cells.Item(18, 5).Value = "  +1.59%  "
$ws.Cells.Item(19, 2).Value = "ShibaInu"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(19, 4).Value = "'0.0000108"
$ws.Cells.Item(19, 5).Value = "  +2.04%  "
$ws.Cells.Item(20, 2).Value = "Uniswap"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(20, 4).Value = "'7.31"
$ws.Cells.Item(20, 5).Value = "  -0.39%  "
$ws.Cells.Item(21, 4).Value = "'75.26"
$ws.Cells.Item(21, 5).Value = "  +1.85%  "
$ws.Cells.Item(22, 4).Value = "'3.38"
$ws.Cells.Item(22, 5).Value = "  -2.89%  "
$ws.Cells.Item(23, 5).Value = "  +9.34%  "
$ws.Cells.Item(24, 4).Value = "'251.84"
$ws.Cells.Item(24, 5).Value = "  +9.66%  "
$ws.Cells.Item(25, 4).Value = "'8.90"
$ws.Cells.Item(25, 5).Value = "  -4.96%  "
$ws.Cells.Item(26, 4).Value = "'11.96"
$ws.Cells.Item(26, 5).Value = "  +1.98%  "
$ws.Cells.Item(27, 5).Value = "  +0.12%  "
$ws.Cells.Item(28, 5).Value = "  +1.06%  "
$ws.Cells.Item(29, 4).Value = "'38.47"
$ws.Cells.Item(29, 5).Value = "  -0.32%  "
$ws.Cells.Item(30, 4).Value = "'22.55"
$ws.Cells.Item(30, 5).Value = "  +6.76%  "
$ws.Cells.Item(31, 4).Value = "'172.89"
$ws.Cells.Item(31, 5).Value = "  -0.84%  "
$ws.Cells.Item(33, 4).Value = "'0.0916"
$ws.Cells.Item(33, 5).Value = "  +3.23%  "
$ws.Cells.Item(34, 4).Value = "'5.90"
$ws.Cells.Item(34, 5).Value = "  +3.40%  "
$ws.Cells.Item(35, 5).Value = "  +3.66%  "
$ws.Cells.Item(36, 4).Value = "'4.92"
$ws.Cells.Item(36, 5).Value = "  +0.29%  "
$ws.Cells.Item(37, 4).Value = "'0.0376"
$ws.Cells.Item(37, 5).Value = "  +1.87%  "
$ws.Cells.Item(38, 4).Value = "'4.07"
$ws.Cells.Item(38, 5).Value = "  -4.71%  "
$ws.Cells.Item(39, 4).Value = "'0.103"
$ws.Cells.Item(39, 5).Value = "  +0.12%  "
$ws.Cells.Item(40, 4).Value = "'2.74"
$ws.Cells.Item(40, 5).Value = "  +10.28%  "
$ws.Cells.Item(41, 4).Value = "'1.50"
$ws.Cells.Item(41, 5).Value = "  +13.80%  "
$ws.Cells.Item(42, 4).Value = "'71.70"
$ws.Cells.Item(42, 5).Value = "  +1.02%  "
$ws.Cells.Item(43, 4).Value = "'0.230"
$ws.Cells.Item(43, 5).Value = "  -2.69%  "
$ws.Cells.Item(44, 5).Value = "  +0.01%  "
$ws.Cells.Item(45, 4).Value = "'12.29"
$ws.Cells.Item(45, 5).Value = "  -4.26%  "
$ws.Cells.Item(46, 4).Value = "'5.62"
$ws.Cells.Item(46, 5).Value = "  +1.89%  "
$ws.Cells.Item(47, 4).Value = "'9.33"
$ws.Cells.Item(47, 5).Value = "  +9.48%  "
$ws.Cells.Item(48, 4).Value = "'110.32"
$ws.Cells.Item(48, 5).Value = "  +5.71%  "
$ws.Cells.Item(49, 5).Value = "  -2.82%  "
$ws.Cells.Item(50, 5).Value = "  +0.94%  "
$ws.Cells.Item(51, 4).Value = "1.490.81"
$ws.Cells.Item(51, 5).Value = "  +4.26%  "
